$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark bug #283 (row 15) as resolved with a "?" in the Done column (C)
$ws.Range("C15").Value = "?"

# Append new bug report as row 23
$ws.Range("A23").Value = -1
$ws.Range("B23").Value = "You can buy dev cards even after there are none left"

# Column B (bug descriptions) should also be centered, in addition to
# already wrapping text
$ws.Range("B1:B23").HorizontalAlignment = -4108
$ws.Range("B1:B23").VerticalAlignment = -4108

# Scroll the view down to show the newly added row, matching where the
# author's cursor ended up after adding the entry
$ws.Application.ActiveWindow.ScrollRow = 17
$ws.Range("A24").Select()
